$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 22
$ws.Cells.Item($row, 1).Value = "Wil je kijken of die bestelling van Van Rijn al is verzonden?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #12: Wil je kijken of die bestelling van Van Rijn al is verzonden?"
$ws.Cells.Item($row, 4).Value = "Bestelling / Levering"
$ws.Cells.Item($row, 5).Value = "Geachte klant,`nDank u voor uw e-mail. Om u beter van dienst te kunnen zijn, zouden wij graag wat meer informatie ontvangen, zoals uw bestelnummer of de naam waaronder de bestelling geplaatst is. Met deze gegevens kunnen wij controleren of de bestelling van Van Rijn al verzonden is.`nWij zien uw reactie graag tegemoet.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Cells.Item($row, 6).Value = "2025-07-23 22:40:39"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

# The multi-line text in column E auto-expands the row height; re-fit it so
# the row keeps the default height (matching the other, unsized rows).
$ws.Rows.Item($row).AutoFit()

# Expand conditional formatting ranges to include the new row, preserving
# the existing rules (priority/dxfId) by moving their AppliesTo range.
$ranges = @("D2:D21", "G2:G21", "H2:H21", "I2:I21", "J2:J21")
foreach ($r in $ranges) {
    $col = $r.Substring(0, 1)
    $fcs = $ws.Range($r).FormatConditions
    $newRange = $ws.Range("$col" + "2:" + "$col" + "22")
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# Update Dashboard summary count for "Bestelling / Levering"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 3
